# GED - Science - Scientific Method Resources Added - Commented Meta Descriptions
# for Exercise Pages
#
# Adds a new "Sheet2" (Watch Hour tracker) after "Sheet1" and makes it the
# active/selected sheet.

$wb = $excel.ActiveWorkbook
$sheet1 = $wb.Worksheets.Item(1)

# Insert the new worksheet right after Sheet1 (becomes "Sheet2" and active).
$ws = $wb.Worksheets.Add($null, $sheet1)

# Column widths.
$ws.Columns.Item(1).ColumnWidth = 12.666666666666666
$ws.Columns.Item(2).ColumnWidth = 12.5
$ws.Columns.Item(3).ColumnWidth = 16.833333333333336

# Header row - write in B, A, C order so shared-string indices line up
# with the source workbook (113=Watch Hour, 114=Date, 115=Daily Average).
$ws.Range("B1").Value = "Watch Hour"
$ws.Range("A1").Value = "Date"
$ws.Range("C1").Value = "Daily Average"

# Date column (formatted d-mmm-yy => numFmtId 15) for rows 2-14.
$dates = 45564, 45565, 45566, 45567, 45568, 45569, 45570, 45571, 45572, 45573, 45574, 45575, 45576
for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 2 + $i
    $cell = $ws.Range("A$row")
    $cell.Value = $dates[$i]
    $cell.NumberFormat = "d-mmm-yy"
}

# Only the first data row has a Watch Hour entry.
$ws.Range("B2").Value = 2761

# Daily Average formulas: C2 is a standalone formula, C3:C14 share one formula.
$ws.Range("C2").Formula = "=B3-B2"
$ws.Range("C3:C14").Formula = "=B4-B3"

# Selection / active cell on the new sheet.
$ws.Range("C3").Select() | Out-Null
